$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date column C for all existing data
#    rows (2..358) from 45182 to 45184.
for ($r = 2; $r -le 358; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# 2) Row 358 picks up an explicit row height (15pt, custom) in the new file.
$ws.Rows.Item(358).RowHeight = 15

# 3) Append a brand-new row (359) with a new logging notification.
$newRow = 359
$ws.Cells.Item($newRow, 1).Value = "A 43172-2023"

$ws.Cells.Item($newRow, 2).Value = 45183
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45184
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item($newRow, 5).Value = "BJURHOLM"
$ws.Cells.Item($newRow, 6).Value = "Holmen skog AB"

$ws.Cells.Item($newRow, 7).Value = 1.1
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Column R keeps the wrap-text style used throughout the sheet, left empty.
$ws.Cells.Item($newRow, 18).Value = ""
$ws.Cells.Item($newRow, 18).WrapText = $true
